# Insert a new data row at row 114, shifting existing rows 114:218 down to 115:219,
# then populate the new row 114 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("114").Insert()

$ws.Range("A114").Value = 10
$ws.Range("B114").Value = "Vega Modelo de Temuco"
$ws.Range("C114").Value = "La Araucanía"
$ws.Range("D114").Value = 44907
$ws.Range("E114").Value = 9
$ws.Range("F114").Value = 100112012
$ws.Range("G114").Value = "Espinaca"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 60
$ws.Range("K114").Value = 9000
$ws.Range("L114").Value = 10000
$ws.Range("M114").Value = 9500
$ws.Range("N114").Value = "$/docena de atados"
$ws.Range("O114").Value = "Región de La Araucanía"
$ws.Range("P114").Value = 3167
$ws.Range("Q114").Value = 3
$ws.Range("R114").Value = "Hortaliza"
